$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.658923
$ws.Range("H2").Value = 7.976769
$ws.Range("I2").Value = 0.01800502032966059
$ws.Range("J2").Value = 0.01800502032966059
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 97.21381066666667
$ws.Range("N2").Value = 291.641432
$ws.Range("O2").Value = 0.2059416914200694
$ws.Range("P2").Value = 0.2059416914200694
$ws.Range("Q2").Value = 258.4840370992454
$ws.Range("R2").Value = 2326.356333893208
$ws.Range("S2").Value = 0.003707984340743037
$ws.Range("T2").Value = 0.003707984340743037

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.658923
$ws.Range("H3").Value = 7.976769
$ws.Range("I3").Value = 0.01800502032966059
$ws.Range("J3").Value = 0.01800502032966059
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 229.1476846666667
$ws.Range("N3").Value = 687.4430540000001
$ws.Range("O3").Value = 0.4854357775055023
$ws.Range("P3").Value = 0.4854357775055023
$ws.Range("Q3").Value = 609.2860491569475
$ws.Range("R3").Value = 5483.574442412527
$ws.Range("S3").Value = 0.008740281042731163
$ws.Range("T3").Value = 0.008740281042731165

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.658923
$ws.Range("H4").Value = 7.976769
$ws.Range("I4").Value = 0.01800502032966059
$ws.Range("J4").Value = 0.01800502032966059
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 145.68382
$ws.Range("N4").Value = 437.05146
$ws.Range("O4").Value = 0.3086225310744283
$ws.Range("P4").Value = 0.3086225310744283
$ws.Range("Q4").Value = 387.36205972586
$ws.Range("R4").Value = 3486.25853753274
$ws.Range("S4").Value = 0.005556754946186389
$ws.Range("T4").Value = 0.005556754946186389

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 122.6832073333333
$ws.Range("H5").Value = 368.049622
$ws.Range("I5").Value = 0.8307550245511554
$ws.Range("J5").Value = 0.8307550245511555
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 97.21381066666667
$ws.Range("N5").Value = 291.641432
$ws.Range("O5").Value = 0.2059416914200694
$ws.Range("P5").Value = 0.2059416914200694
$ws.Range("Q5").Value = 11926.50208968208
$ws.Range("R5").Value = 107338.5188071387
$ws.Range("S5").Value = 0.1710870949117862
$ws.Range("T5").Value = 0.1710870949117862

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 122.6832073333333
$ws.Range("H6").Value = 368.049622
$ws.Range("I6").Value = 0.8307550245511554
$ws.Range("J6").Value = 0.8307550245511555
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 229.1476846666667
$ws.Range("N6").Value = 687.4430540000001
$ws.Range("O6").Value = 0.4854357775055023
$ws.Range("P6").Value = 0.4854357775055023
$ws.Range("Q6").Value = 28112.57290791396
$ws.Range("R6").Value = 253013.1561712256
$ws.Range("S6").Value = 0.4032782112595927
$ws.Range("T6").Value = 0.4032782112595928

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 122.6832073333333
$ws.Range("H7").Value = 368.049622
$ws.Range("I7").Value = 0.8307550245511554
$ws.Range("J7").Value = 0.8307550245511555
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 145.68382
$ws.Range("N7").Value = 437.05146
$ws.Range("O7").Value = 0.3086225310744283
$ws.Range("P7").Value = 0.3086225310744283
$ws.Range("Q7").Value = 17872.95829417201
$ws.Range("R7").Value = 160856.6246475481
$ws.Range("S7").Value = 0.2563897183797764
$ws.Range("T7").Value = 0.2563897183797764

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 22.334626
$ws.Range("H8").Value = 67.003878
$ws.Range("I8").Value = 0.1512399551191839
$ws.Range("J8").Value = 0.151239955119184
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 97.21381066666667
$ws.Range("N8").Value = 291.641432
$ws.Range("O8").Value = 0.2059416914200694
$ws.Range("P8").Value = 0.2059416914200694
$ws.Range("Q8").Value = 2171.234103274811
$ws.Range("R8").Value = 19541.1069294733
$ws.Range("S8").Value = 0.03114661216754012
$ws.Range("T8").Value = 0.03114661216754012

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 22.334626
$ws.Range("H9").Value = 67.003878
$ws.Range("I9").Value = 0.1512399551191839
$ws.Range("J9").Value = 0.151239955119184
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 229.1476846666667
$ws.Range("N9").Value = 687.4430540000001
$ws.Range("O9").Value = 0.4854357775055023
$ws.Range("P9").Value = 0.4854357775055023
$ws.Range("Q9").Value = 5117.927835795936
$ws.Range("R9").Value = 46061.35052216342
$ws.Range("S9").Value = 0.07341728520317833
$ws.Range("T9").Value = 0.07341728520317835

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 22.334626
$ws.Range("H10").Value = 67.003878
$ws.Range("I10").Value = 0.1512399551191839
$ws.Range("J10").Value = 0.151239955119184
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 145.68382
$ws.Range("N10").Value = 437.05146
$ws.Range("O10").Value = 0.3086225310744283
$ws.Range("P10").Value = 0.3086225310744283
$ws.Range("Q10").Value = 3253.79363395132
$ws.Range("R10").Value = 29284.14270556188
$ws.Range("S10").Value = 0.0466760577484655
$ws.Range("T10").Value = 0.0466760577484655

